$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet stores Price/Volume columns as plain text (inline strings) in the
# source workbook, preserving formats like trailing zeros ("7.00"), thousand-dot
# grouping ("61.571.08") and small-number notation ("0.0000107"). Force the cells
# to Text format before assigning so Excel does not silently reinterpret the
# strings as numbers (which would strip trailing zeros / use scientific notation).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.571.08"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.46%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.893.87"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.00%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.48%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.57"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.24%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.11%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.505"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.891.93"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.99%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.85%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.56%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.431"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.21%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.54%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.87"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.78%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.58%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.374.97"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.96%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.587.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.32%  "

# Row 18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.97%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.904.67"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.78%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.87"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.18%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.06"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.15%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.656"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.76%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.81"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.88%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.28"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.11%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.98"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.76%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.99"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -10.55%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.01"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.59%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000107"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.48%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.01"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.80%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -7.25%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.06%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.106"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.96%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.49"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.65%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.958"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.50%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.39"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.82"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.71%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.93"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.54%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.82"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -8.88%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.21"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.55%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.114"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.66%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.48"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.94%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.266"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.36%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.696.66"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.19%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.52"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.22%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.56%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "347.50"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.18%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.45%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.57"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -5.40%  "
